$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1853146853146853
$ws.Range("C2").Value = 0.5454545454545454
$ws.Range("J2").Value = 0.01048951048951049
$ws.Range("P2").Value = 0.1398601398601399
$ws.Range("S2").Value = 0.1188811188811189

$ws.Range("C3").Value = 0.03144654088050314
$ws.Range("J3").Value = 0.01257861635220126
$ws.Range("P3").Value = 0.7295597484276729
$ws.Range("S3").Value = 0.2264150943396226

$ws.Range("J4").Value = 0.04081632653061224
$ws.Range("P4").Value = 0.6326530612244898
$ws.Range("S4").Value = 0.3265306122448979

$ws.Range("B6").Value = 0.06493506493506493
$ws.Range("D6").Value = 0.008658008658008658
$ws.Range("F6").Value = 0.09956709956709957
$ws.Range("J6").Value = 0.2121212121212121
$ws.Range("O6").Value = 0.0303030303030303
$ws.Range("Q6").Value = 0.1861471861471861
$ws.Range("R6").Value = 0.07792207792207792
$ws.Range("S6").Value = 0.3203463203463203

$ws.Range("B7").Value = 0.1218274111675127
$ws.Range("D7").Value = 0.03045685279187817
$ws.Range("F7").Value = 0.07614213197969544
$ws.Range("J7").Value = 0.1065989847715736
$ws.Range("O7").Value = 0.02538071065989848
$ws.Range("Q7").Value = 0.2182741116751269
$ws.Range("R7").Value = 0.05076142131979695
$ws.Range("S7").Value = 0.3705583756345178

$ws.Range("B8").Value = 0.1070559610705596
$ws.Range("D8").Value = 0.0194647201946472
$ws.Range("F8").Value = 0.05596107055961071
$ws.Range("J8").Value = 0.09245742092457421
$ws.Range("O8").Value = 0.0267639902676399
$ws.Range("Q8").Value = 0.218978102189781
$ws.Range("R8").Value = 0.06569343065693431
$ws.Range("S8").Value = 0.413625304136253

$ws.Range("B9").Value = 0.1323529411764706
$ws.Range("F9").Value = 0.1274509803921569
$ws.Range("J9").Value = 0.09313725490196079
$ws.Range("O9").Value = 0.004901960784313725
$ws.Range("Q9").Value = 0.196078431372549
$ws.Range("R9").Value = 0.08333333333333333
$ws.Range("S9").Value = 0.3333333333333333

$ws.Range("B10").Value = 0.1034790365744871
$ws.Range("D10").Value = 0.0312221231043711
$ws.Range("E10").Value = 0.0008920606601248885
$ws.Range("F10").Value = 0.07582515611061552
$ws.Range("J10").Value = 0.1034790365744871
$ws.Range("O10").Value = 0.0231935771632471
$ws.Range("Q10").Value = 0.2185548617305977
$ws.Range("R10").Value = 0.05798394290811775
$ws.Range("S10").Value = 0.3853702051739518

$ws.Range("G11").Value = 0.1007194244604317
$ws.Range("J11").Value = 0.06474820143884892
$ws.Range("K11").Value = 0.1690647482014389
$ws.Range("L11").Value = 0.6510791366906474
$ws.Range("S11").Value = 0.01438848920863309

$ws.Range("G12").Value = 0.78125
$ws.Range("J12").Value = 0.1458333333333333
$ws.Range("L12").Value = 0.04166666666666666
$ws.Range("S12").Value = 0.03125

$ws.Range("G13").Value = 0.5714285714285714
$ws.Range("J13").Value = 0.3714285714285714
$ws.Range("S13").Value = 0.05714285714285714

$ws.Range("F15").Value = 0.01746724890829694
$ws.Range("H15").Value = 0.1834061135371179
$ws.Range("I15").Value = 0.07860262008733625
$ws.Range("J15").Value = 0.3100436681222707
$ws.Range("K15").Value = 0.04803493449781659
$ws.Range("M15").Value = 0.004366812227074236
$ws.Range("O15").Value = 0.0611353711790393
$ws.Range("S15").Value = 0.2969432314410481

$ws.Range("F16").Value = 0.02209944751381215
$ws.Range("H16").Value = 0.1602209944751381
$ws.Range("I16").Value = 0.06629834254143646
$ws.Range("J16").Value = 0.4088397790055249
$ws.Range("K16").Value = 0.138121546961326
$ws.Range("M16").Value = 0.005524861878453038
$ws.Range("N16").Value = 0.005524861878453038
$ws.Range("O16").Value = 0.08287292817679558
$ws.Range("S16").Value = 0.1104972375690608

$ws.Range("F17").Value = 0.0131578947368421
$ws.Range("H17").Value = 0.1776315789473684
$ws.Range("I17").Value = 0.1140350877192982
$ws.Range("J17").Value = 0.4210526315789473
$ws.Range("K17").Value = 0.1074561403508772
$ws.Range("M17").Value = 0.01754385964912281
$ws.Range("N17").Value = 0.008771929824561403
$ws.Range("O17").Value = 0.05701754385964912
$ws.Range("S17").Value = 0.08333333333333333

$ws.Range("F18").Value = 0.007352941176470588
$ws.Range("H18").Value = 0.1838235294117647
$ws.Range("I18").Value = 0.08088235294117647
$ws.Range("J18").Value = 0.4191176470588235
$ws.Range("K18").Value = 0.1029411764705882
$ws.Range("O18").Value = 0.1029411764705882
$ws.Range("S18").Value = 0.1029411764705882

$ws.Range("F19").Value = 0.01824500434404865
$ws.Range("H19").Value = 0.2067767158992181
$ws.Range("I19").Value = 0.09643788010425716
$ws.Range("J19").Value = 0.3753258036490009
$ws.Range("K19").Value = 0.1129452649869679
$ws.Range("M19").Value = 0.02258905299739357
$ws.Range("N19").Value = 0.0008688097306689834
$ws.Range("O19").Value = 0.0686359687228497
$ws.Range("S19").Value = 0.09817549956559514
